# Issue #840 Spring IO Platform translation
#
# The callout textbox on slide 2 that used to just read "parent" (labelling
# the parent pom/module relationship) is widened, shifted, switched to
# word-wrapping text, and reworded to describe the dependency-management
# import step.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(17)

# Reposition / resize the textbox (values are EMU/12700 expressed in points;
# nudged by half an EMU-in-points so the point->EMU round trip lands on the
# exact target EMU instead of being truncated one EMU short).
$shp.Left   = 46.00011826098815
$shp.Top    = 330.54066467922024
$shp.Width  = 333.8378295935506
$shp.Height = 29.081298833350314

# Allow the now-wider text to wrap instead of forcing a single line.
$shp.TextFrame.WordWrap = $true

# Update the label text.
$shp.TextFrame.TextRange.Text = "Import the dependency management"
